# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph
# together with the blank paragraph before it, the blank paragraph
# after it, and the blank page-break paragraph that followed it.
# This collapses the four paragraphs that sat between the
# "LOB1036: Geometria Analítica (Requisito fraco)" paragraph and the
# trailing blank / page-break paragraphs at the end of the document.

$d = $word.ActiveDocument

$target = $d.Content.Find
$found = $target.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $hit = $d.Range($target.Parent.Start, $target.Parent.End)
    $hitParaIndex = $hit.Paragraphs(1).Index

    $startPara = $hitParaIndex - 1
    $endPara = $hitParaIndex + 2

    $start = $d.Paragraphs($startPara).Range.Start
    $end = $d.Paragraphs($endPara).Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
